# Update "Förändrad" (changed) date column C for rows 2-6 from 2023-09-14 (45183)
# to 2023-09-15 (45184), matching the automatic update of files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
